$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.3464964993005633, 1.65323645889881, 2938.103010863317, 71517.89157740913, 0, 74457.99432123064)
    3 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    5 = @(0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.719096746035642)
    6 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    7 = @(0.06328177979961902, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 1, 8.35085209798723)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
